$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 8 label changes from "Model" to "production_function"
$ws.Range("A8").Value = "production_function"

# Insert a new row 9 ("L_curve") above the old row 9, shifting rows 9-17 down to 10-18
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0

# Give the two new label cells (A8/A9) the same style as the other section-header
# label cell on this sheet (A1), which is cellXf index 3 (bold-ish Verdana header style)
$ws.Range("A8:A9").Style = $ws.Range("A1").Style

# Restore exact numeric literals on the two rows that shifted down, in case the
# structural shift re-serialized their floating point text representation
$ws.Range("E14").Value = 1.6
$ws.Range("R18").Value = 1.6

# Make this sheet the active tab (was threshold_b, now optimization_parameters)
$ws.Activate()
$ws.Range("A9:B9").Select()
